$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capital cost table (rows 3-17): refreshed simulation results ---
$ws.Range("C3").Value  = 65.56319540899871
$ws.Range("C4").Value  = 32.40560529144166
$ws.Range("C5").Value  = 2.622527816359948
$ws.Range("C6").Value  = 5.900687586809884
$ws.Range("C7").Value  = 2.950343793404942
$ws.Range("C8").Value  = 109.4423598970151
$ws.Range("C9").Value  = 10.94423598970151
$ws.Range("C10").Value = 10.94423598970151
$ws.Range("C11").Value = 21.88847197940303
$ws.Range("C12").Value = 32.83270796910454
$ws.Range("C13").Value = 10.94423598970151
$ws.Range("C14").Value = 87.55388791761214
$ws.Range("C15").Value = 196.9962478146273
$ws.Range("C16").Value = 9.849812390731365
$ws.Range("C17").Value = 206.8460602053586

# --- Raw materials / by-products table (rows 21-29) ---
# Unmerge the old A21:A22 and A24:A28 blocks before rearranging the rows.
$ws.Range("A21:A22").UnMerge()
$ws.Range("A24:A28").UnMerge()

# Row 21 now holds the "By-products and credits" section (just Wastewater).
$ws.Range("A21").Value = "By-products and credits"
$ws.Range("B21").Value = "Wastewater"
$ws.Range("C21").Value = -1.962995609856692
$ws.Range("D21").Value = -6.266069596713132

# Row 22 starts the "Raw materials" section (merged label through row 28).
$ws.Range("A22").Value = "Raw materials"
$ws.Range("B22").Value = "Glucose"
$ws.Range("C22").Value = 240.404025
$ws.Range("D22").Value = 144.7370872673688

$ws.Range("A23").Value = ""
$ws.Range("B23").Value = "Process water"
$ws.Range("C23").Value = 0.320236305
$ws.Range("D23").Value = 0.9954177967561536

$ws.Range("A24").Value = ""
$ws.Range("B24").Value = "Tridecane"
$ws.Range("C24").Value = 878.1550799999999
$ws.Range("D24").Value = 0.002402636063926985

$ws.Range("A25").Value = ""
$ws.Range("B25").Value = "CSL"
$ws.Range("C25").Value = 51.528108
$ws.Range("D25").Value = 0.5665456791641987

$ws.Range("A26").Value = ""
$ws.Range("B26").Value = "DAP"
$ws.Range("C26").Value = 895.3915949999999
$ws.Range("D26").Value = 1.181188778425073

$ws.Range("A27").Value = ""
$ws.Range("B27").Value = "Salt"
$ws.Range("C27").Value = 136.07775
$ws.Range("D27").Value = 4.230790853725686

$ws.Range("A28").Value = ""
$ws.Range("B28").Value = "Natural gas"
$ws.Range("C28").Value = 197.76633
$ws.Range("D28").Value = 3.954929396805927

$ws.Range("D29").Value = 161.934432005023

# Re-merge the "Raw materials" label across its new row span.
$ws.Range("A22:A28").Merge()

# --- Labor/overhead table (rows 35-36): refreshed simulation results ---
$ws.Range("C35").Value = 1.966895862269961
$ws.Range("D35").Value = 1.888220027779163
$ws.Range("C36").Value = 0.458942367862991
$ws.Range("D36").Value = 0.4405846731484713
